# Generate Report for handoff
# Marks the e900ff5d-9406-4620-b0b0-4ecc073d7efd.md file as "Ready for
# handoff" on the Overview sheet and on each language sheet, and records
# the new "Latest Handoff Datetime" for that handoff on each language
# sheet (the "Latest Handback DateTime" is left untouched).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-17 07:35:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-17 07:35:35"
